$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 852.7143
$ws.Cells.Item(41, 9).Value = 1001.3125
$ws.Cells.Item(41, 10).Value = 377.2
$ws.Cells.Item(41, 11).Value = 1001.3125
$ws.Cells.Item(41, 12).Value = 377.2
$ws.Cells.Item(41, 13).Value = -561.3125
$ws.Cells.Item(41, 14).Value = -1257.2

$ws.Cells.Item(116, 8).Value = 3287.5
$ws.Cells.Item(116, 9).Value = 2625.2
$ws.Cells.Item(116, 11).Value = 2625.2
$ws.Cells.Item(116, 13).Value = 816.8000000000002

$ws.Cells.Item(125, 8).Value = 1138.6111
$ws.Cells.Item(125, 9).Value = 992.7857
$ws.Cells.Item(125, 11).Value = 8935.0713
$ws.Cells.Item(125, 13).Value = -6475.0713

$ws.Cells.Item(129, 8).Value = 1999.9584
$ws.Cells.Item(129, 9).Value = 746.7857
$ws.Cells.Item(129, 11).Value = 2240.3571
$ws.Cells.Item(129, 13).Value = 2759.6429

$ws.Cells.Item(132, 8).Value = 3601.543
$ws.Cells.Item(132, 9).Value = 2045.2424
$ws.Cells.Item(132, 11).Value = 6135.7272
$ws.Cells.Item(132, 13).Value = -3605.7272

$ws.Cells.Item(137, 8).Value = 1341.6923
$ws.Cells.Item(137, 9).Value = 1194.9
$ws.Cells.Item(137, 10).Value = 1831
$ws.Cells.Item(137, 11).Value = 3584.7
$ws.Cells.Item(137, 12).Value = 5493
$ws.Cells.Item(137, 13).Value = -1034.7
$ws.Cells.Item(137, 14).Value = -10593

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 945.8542
$ws.Cells.Item(2, 9).Value = 879.6977000000001
$ws.Cells.Item(2, 11).Value = 879.6977000000001
$ws.Cells.Item(2, 13).Value = -766.6977000000001

$ws.Cells.Item(31, 8).Value = 6847.5835
$ws.Cells.Item(31, 9).Value = 4695.636
$ws.Cells.Item(31, 11).Value = 4695.636
$ws.Cells.Item(31, 13).Value = -4401.636

$ws.Cells.Item(61, 8).Value = 1480.4286
$ws.Cells.Item(61, 9).Value = 1480.4286
$ws.Cells.Item(61, 11).Value = 1480.4286
$ws.Cells.Item(61, 13).Value = -1268.4286

$ws.Cells.Item(70, 8).Value = 42786
$ws.Cells.Item(70, 10).Value = 42786
$ws.Cells.Item(70, 12).Value = 42786
$ws.Cells.Item(70, 14).Value = -43326

$ws.Cells.Item(73, 8).Value = 42786
$ws.Cells.Item(73, 10).Value = 42786
$ws.Cells.Item(73, 12).Value = 42786
$ws.Cells.Item(73, 14).Value = -44658

$ws.Cells.Item(74, 8).Value = 2163.5881
$ws.Cells.Item(74, 9).Value = 1727.4166
$ws.Cells.Item(74, 11).Value = 1727.4166
$ws.Cells.Item(74, 13).Value = -853.4166

$ws.Cells.Item(77, 8).Value = 2163.5881
$ws.Cells.Item(77, 9).Value = 1727.4166
$ws.Cells.Item(77, 11).Value = 8637.083000000001
$ws.Cells.Item(77, 13).Value = -4269.083000000001

$ws.Cells.Item(116, 8).Value = 945.8542
$ws.Cells.Item(116, 9).Value = 879.6977000000001
$ws.Cells.Item(116, 11).Value = 879.6977000000001
$ws.Cells.Item(116, 13).Value = 1414.3023

$ws.Cells.Item(122, 8).Value = 1981.7805
$ws.Cells.Item(122, 9).Value = 1449.138
$ws.Cells.Item(122, 11).Value = 4347.414
$ws.Cells.Item(122, 13).Value = -1897.414

$ws.Cells.Item(132, 8).Value = 3488.8
$ws.Cells.Item(132, 9).Value = 2913.2942
$ws.Cells.Item(132, 11).Value = 8739.882599999999
$ws.Cells.Item(132, 13).Value = -6209.882599999999

$ws.Cells.Item(136, 8).Value = 1480.4286
$ws.Cells.Item(136, 9).Value = 1480.4286
$ws.Cells.Item(136, 11).Value = 4441.2858
$ws.Cells.Item(136, 13).Value = -1891.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 945.8542
$ws.Cells.Item(3, 9).Value = 879.6977000000001
$ws.Cells.Item(3, 11).Value = 879.6977000000001
$ws.Cells.Item(3, 13).Value = -765.6977000000001

$ws.Cells.Item(94, 8).Value = 1104.5
$ws.Cells.Item(94, 9).Value = 1200
$ws.Cells.Item(94, 11).Value = 1200
$ws.Cells.Item(94, 13).Value = -749

$ws.Cells.Item(99, 8).Value = 7780.5806
$ws.Cells.Item(99, 9).Value = 8217.173000000001
$ws.Cells.Item(99, 11).Value = 8217.173000000001
$ws.Cells.Item(99, 13).Value = -6719.173000000001

$ws.Cells.Item(107, 8).Value = 88003
$ws.Cells.Item(107, 9).Value = 250382
$ws.Cells.Item(107, 10).Value = 6813.5
$ws.Cells.Item(107, 11).Value = 250382
$ws.Cells.Item(107, 12).Value = 6813.5
$ws.Cells.Item(107, 13).Value = -248462
$ws.Cells.Item(107, 14).Value = -10653.5

$ws.Cells.Item(134, 8).Value = 992.4595
$ws.Cells.Item(134, 9).Value = 737.3823
$ws.Cells.Item(134, 11).Value = 2212.1469
$ws.Cells.Item(134, 13).Value = 322.8531000000003

$ws.Cells.Item(135, 8).Value = 47500
$ws.Cells.Item(135, 10).Value = 47500
$ws.Cells.Item(135, 12).Value = 47500
$ws.Cells.Item(135, 14).Value = -57640

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 839.5333000000001
$ws.Cells.Item(22, 9).Value = 440.16666
$ws.Cells.Item(22, 10).Value = 1105.7778
$ws.Cells.Item(22, 11).Value = 440.16666
$ws.Cells.Item(22, 12).Value = 1105.7778
$ws.Cells.Item(22, 13).Value = -90.16665999999998
$ws.Cells.Item(22, 14).Value = -1805.7778

$ws.Cells.Item(31, 8).Value = 3897.35
$ws.Cells.Item(31, 9).Value = 2133.0833
$ws.Cells.Item(31, 11).Value = 2133.0833
$ws.Cells.Item(31, 13).Value = -1838.0833

$ws.Cells.Item(34, 8).Value = 3897.35
$ws.Cells.Item(34, 9).Value = 2133.0833
$ws.Cells.Item(34, 11).Value = 2133.0833
$ws.Cells.Item(34, 13).Value = -1931.0833

$ws.Cells.Item(52, 8).Value = 101250
$ws.Cells.Item(52, 10).Value = 124386.5
$ws.Cells.Item(52, 12).Value = 124386.5
$ws.Cells.Item(52, 14).Value = -124974.5

$ws.Cells.Item(112, 8).Value = 51661.668
$ws.Cells.Item(112, 10).Value = 51661.668
$ws.Cells.Item(112, 12).Value = 51661.668
$ws.Cells.Item(112, 14).Value = -54615.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 201.22223
$ws.Cells.Item(40, 9).Value = 143.66667
$ws.Cells.Item(40, 10).Value = 316.33334
$ws.Cells.Item(40, 11).Value = 574.66668
$ws.Cells.Item(40, 12).Value = 1265.33336
$ws.Cells.Item(40, 13).Value = -505.66668
$ws.Cells.Item(40, 14).Value = -1403.33336

$ws.Cells.Item(88, 8).Value = 11247.75
$ws.Cells.Item(88, 10).Value = 11247.75
$ws.Cells.Item(88, 12).Value = 33743.25
$ws.Cells.Item(88, 14).Value = -34599.25

$ws.Cells.Item(91, 8).Value = 11247.75
$ws.Cells.Item(91, 10).Value = 11247.75
$ws.Cells.Item(91, 12).Value = 33743.25
$ws.Cells.Item(91, 14).Value = -36707.25

$ws.Cells.Item(122, 8).Value = 853.1111
$ws.Cells.Item(122, 9).Value = 594.5
$ws.Cells.Item(122, 11).Value = 5350.5
$ws.Cells.Item(122, 13).Value = -2900.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 80429.336
$ws.Cells.Item(32, 10).Value = 75644
$ws.Cells.Item(32, 12).Value = 75644
$ws.Cells.Item(32, 14).Value = -76236

$ws.Cells.Item(101, 8).Value = 40157
$ws.Cells.Item(101, 10).Value = 40157
$ws.Cells.Item(101, 12).Value = 40157
$ws.Cells.Item(101, 14).Value = -46647

$ws.Cells.Item(107, 8).Value = 29414762
$ws.Cells.Item(107, 9).Value = 626.8
$ws.Cells.Item(107, 10).Value = 41670652
$ws.Cells.Item(107, 11).Value = 626.8
$ws.Cells.Item(107, 12).Value = 41670652
$ws.Cells.Item(107, 13).Value = 1293.2
$ws.Cells.Item(107, 14).Value = -41674492

$ws.Cells.Item(122, 8).Value = 2908.5264
$ws.Cells.Item(122, 9).Value = 2447.2856
$ws.Cells.Item(122, 10).Value = 4200
$ws.Cells.Item(122, 11).Value = 7341.8568
$ws.Cells.Item(122, 12).Value = 12600
$ws.Cells.Item(122, 13).Value = -4891.8568
$ws.Cells.Item(122, 14).Value = -17500

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 556.24
$ws.Cells.Item(16, 9).Value = 575.6667
$ws.Cells.Item(16, 11).Value = 575.6667
$ws.Cells.Item(16, 13).Value = -405.6667

$ws.Cells.Item(40, 8).Value = 15714.214
$ws.Cells.Item(40, 9).Value = 25599.8
$ws.Cells.Item(40, 10).Value = 10222.223
$ws.Cells.Item(40, 11).Value = 25599.8
$ws.Cells.Item(40, 12).Value = 10222.223
$ws.Cells.Item(40, 13).Value = -25463.8
$ws.Cells.Item(40, 14).Value = -10494.223

$ws.Cells.Item(93, 8).Value = 16902.32
$ws.Cells.Item(93, 9).Value = 3187.5557
$ws.Cells.Item(93, 10).Value = 52168.855
$ws.Cells.Item(93, 11).Value = 3187.5557
$ws.Cells.Item(93, 12).Value = 52168.855
$ws.Cells.Item(93, 13).Value = -1939.5557
$ws.Cells.Item(93, 14).Value = -54664.855

$ws.Cells.Item(97, 8).Value = 20000
$ws.Cells.Item(97, 10).Value = 20000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 14).Value = -21982

$ws.Cells.Item(122, 8).Value = 94577.82000000001
$ws.Cells.Item(122, 9).Value = 136109.73
$ws.Cells.Item(122, 11).Value = 408329.1900000001
$ws.Cells.Item(122, 13).Value = -405879.1900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 9000
$ws.Cells.Item(18, 9).Value = 9000
$ws.Cells.Item(18, 11).Value = 9000
$ws.Cells.Item(18, 13).Value = -8827

$ws.Cells.Item(96, 8).Value = 1889.1111
$ws.Cells.Item(96, 9).Value = 1783.6666
$ws.Cells.Item(96, 10).Value = 2100
$ws.Cells.Item(96, 11).Value = 1783.6666
$ws.Cells.Item(96, 12).Value = 2100
$ws.Cells.Item(96, 13).Value = -410.6666
$ws.Cells.Item(96, 14).Value = -4846

$ws.Cells.Item(107, 8).Value = 13160087
$ws.Cells.Item(107, 9).Value = 1988.5
$ws.Cells.Item(107, 11).Value = 5965.5
$ws.Cells.Item(107, 13).Value = -4045.5

$ws.Cells.Item(112, 8).Value = 36661.668
$ws.Cells.Item(112, 10).Value = 36661.668
$ws.Cells.Item(112, 12).Value = 36661.668
$ws.Cells.Item(112, 14).Value = -39615.668

$ws.Cells.Item(132, 8).Value = 2163.2112
$ws.Cells.Item(132, 9).Value = 1681.6558
$ws.Cells.Item(132, 11).Value = 5044.9674
$ws.Cells.Item(132, 13).Value = -2514.9674

$ws.Cells.Item(136, 8).Value = 789.2857
$ws.Cells.Item(136, 9).Value = 753.67645
$ws.Cells.Item(136, 11).Value = 2261.02935
$ws.Cells.Item(136, 13).Value = 288.9706499999998
